$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: add date (A29), batch value (B29), CRM bottle note (F29)
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 43504
$ws.Range("B29").Value = 2206.2235349408302
$ws.Range("F29").Value = "New CRM bottle"

# Update selection to G29
$ws.Range("G29").Select()
